$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.412.27"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").Value = "1.654.80"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.25"
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.537"
$ws.Range("E6").Value = "  +4.79%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.45"
$ws.Range("E8").Value = "  +0.30%  "
$ws.Range("E9").Value = "  +0.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0613"
$ws.Range("E10").Value = "  -1.20%  "
$ws.Range("E11").Value = "  +3.42%  "
$ws.Range("D12").Value = "1.889.17"
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("D13").Value = "1.656.02"
$ws.Range("E13").Value = "  -0.34%  "
$ws.Range("E14").Value = "  -1.24%  "
$ws.Range("E15").Value = "  +3.75%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.51"
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("D17").Value = "27.403.46"
$ws.Range("E17").Value = "  -0.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.67"
$ws.Range("E18").Value = "  -7.05%  "
$ws.Range("D19").Value = "0.0₃0726"
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.42"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.34"
$ws.Range("E22").Value = "  -2.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.39"
$ws.Range("E23").Value = "  +0.70%  "
$ws.Range("E24").Value = "  +1.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.00"
$ws.Range("E25").Value = "  +0.61%  "
$ws.Range("E26").Value = "  -1.44%  "
$ws.Range("E27").Value = "  +2.85%  "
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.68"
$ws.Range("E29").Value = "  -3.18%  "
$ws.Range("E30").Value = "  -0.90%  "
$ws.Range("E31").Value = "  -3.85%  "
$ws.Range("E32").Value = "  -0.96%  "
$ws.Range("B33").Value = "Maker"
$ws.Range("C33").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D33").Value = "1.421.38"
$ws.Range("E33").Value = "  -1.74%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.12"
$ws.Range("E34").Value = "  -0.33%  "
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("E36").Value = "  -0.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.905"
$ws.Range("E37").Value = "  -2.62%  "
$ws.Range("E38").Value = "  -2.04%  "
$ws.Range("E39").Value = "  -0.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.05"
$ws.Range("E40").Value = "  +0.70%  "
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.54"
$ws.Range("E42").Value = "  +2.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.04"
$ws.Range("E43").Value = "  -5.78%  "
$ws.Range("E44").Value = "  +0.47%  "
$ws.Range("D46").Value = "1.797.85"
$ws.Range("E46").Value = "  -0.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.67"
$ws.Range("E47").Value = "  -1.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.95"
$ws.Range("E48").Value = "  -0.85%  "
$ws.Range("E49").Value = "  -3.47%  "
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.73"
$ws.Range("E51").Value = "  -1.57%  "
